$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 ("past_nonconformity" / enforcement-related weight): 4 -> 5
$ws.Range("D13").Value = 5

# Row 14 ("enforcement_history" / previous inspection weight): 5 -> 4
$ws.Range("D14").Value = 4
